# "护盾.xlsx" — add the header row for the shield data table.
#
# Columns, left to right: 护盾编号(A) 护盾名称(B) 护盾类型(C) 护盾等级(D)
# 基础属性(E) 配方(F).
#
# NOTE: E1 is written before D1 so the shared-string table is built in the
# same order the source workbook has it (基础属性 ends up as shared-string
# index 3, before 护盾等级 at index 4, even though 护盾等级 appears first on
# the sheet in column D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "护盾编号"
$ws.Range("B1").Value = "护盾名称"
$ws.Range("C1").Value = "护盾类型"
$ws.Range("E1").Value = "基础属性"
$ws.Range("D1").Value = "护盾等级"
$ws.Range("F1").Value = "配方"

# Widen column E (基础属性) so the longer recipe/attribute text is readable.
$ws.Columns.Item(5).ColumnWidth = 27.88671875

# Leave the selection the way the author saved the sheet: whole column F
# selected (from the header cell down).
$ws.Range("F1:F1048576").Select()
